$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the bordered cell style (style index 1) from the existing data rows
# down onto the 12 new rows before filling in values.
$ws.Range("A9:K10").Copy() | Out-Null
$ws.Range("A11:K22").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A11").Value = "Задание 7"
$ws.Range("B11").Value = "Политика 7"
$ws.Range("C11").Value = "Задание 7"
$ws.Range("D11").Value = "promoall.docx"
$ws.Range("D12").Value = "promotwo.docx"
$ws.Range("D13").Value = "promoall.docx"
$ws.Range("D14").Value = "promotwo.docx"
$ws.Range("E11:E12").Value = "Security"
$ws.Range("E13:E14").Value = "Sales"
$ws.Range("F11:F14").Value = "External"
$ws.Range("G11:G12").Value = "Политика 7"
$ws.Range("H11:H14").Value = "Задание 7"
$ws.Range("I11:I12").Value = "Forbidden"
$ws.Range("I13:I14").Value = "Allowed"
$ws.Range("J11:J12").Value = "Medium"
$ws.Range("J13:J14").Value = "No"
$ws.Range("K11:K12").Value = "Политика 7"
$ws.Range("A15").Value = "Задание 8"
$ws.Range("B15").Value = "Политика 8"
$ws.Range("C15").Value = "Задание 8"
$ws.Range("D15").Value = "catoo.jpg"
$ws.Range("D16").Value = "catoo_reduced.jpg"
$ws.Range("D17").Value = "catoo.jpg"
$ws.Range("D18").Value = "catoo_reduced.jpg"
$ws.Range("E15:E18").Value = "IT"
$ws.Range("F15:F16").Value = "External"
$ws.Range("F17:F18").Value = "Security"
$ws.Range("G15").Value = "Политика 8"
$ws.Range("G17").Value = "Политика 8"
$ws.Range("H15").Value = "Задание 8"
$ws.Range("H17").Value = "Задание 8"
$ws.Range("I15").Value = "Forbidden"
$ws.Range("I16").Value = "Allowed"
$ws.Range("I17").Value = "Forbidden"
$ws.Range("I18").Value = "Allowed"
$ws.Range("J15").Value = "Low"
$ws.Range("J16").Value = "No"
$ws.Range("J17").Value = "Low"
$ws.Range("J18").Value = "No"
$ws.Range("K15").Value = "Политика 8"
$ws.Range("K17").Value = "Политика 8"
$ws.Range("A19").Value = "Задание 9"
$ws.Range("B19").Value = "Политика 9"
$ws.Range("C19").Value = "Задание 9"
$ws.Range("D19").Value = "evil_routes.rtf"
$ws.Range("D20").Value = "routes.rtf"
$ws.Range("D21").Value = "evil_routes.rtf"
$ws.Range("D22").Value = "routes.rtf"
$ws.Range("E19:E20").Value = "Sales"
$ws.Range("E21:E22").Value = "Accounting"
$ws.Range("F19:F20").Value = "IT"
$ws.Range("F21:F22").Value = "Sales"
$ws.Range("G19:G20").Value = "Политика 9"
$ws.Range("H19:H22").Value = "Задание 9"
$ws.Range("I19:I22").Value = "Allowed"
$ws.Range("J19:J20").Value = "Medium"
$ws.Range("J21:J22").Value = "No"
$ws.Range("K19:K20").Value = "Политика 9"

# Update the used dimension and the saved view/selection state to match the
# grown data range.
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("D25").Select() | Out-Null
